$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (date serial, y_0, y_0_forecast, y_1, y_1_forecast) for rows 2..19
$data = @(
    @(39400, 2007, 2.070003986395053, 2008, $null),
    @(39765, 2008, 0.517569958955022, 2009, $null),
    @(40130, 2009, -3.956152295564896, 2010, $null),
    @(40494, 2010, 1.234995474941392, 2011, 0.7327527981400461),
    @(40862, 2011, 0.899360810820804, 2012, 0.4113848771853501),
    @(41228, 2012, 0.9010266119894084, 2013, 1.203907967581297),
    @(41592, 2013, 0.02019328874804938, 2014, -0.03860754389363175),
    @(41957, 2014, 0.1729981757035093, 2015, 0.2629870913912535),
    @(42321, 2015, 0.09752710595589686, 2016, 0.1987429576382649),
    @(42689, 2016, -0.5280591151586633, 2017, 0.1903092973221776),
    @(43053, 2017, 0.07201851318385799, 2018, 0.2843016498274009),
    @(43418, 2018, 0.3727661260635617, 2019, -0.8612142616933327),
    @(43783, 2019, -0.801759526476209, 2020, 0.06491682578968483),
    @(44159, 2020, -1.103489789942047, 2021, 1.323658311025055),
    @(44525, 2021, 0.9704846793491928, 2022, -0.6989646400249128),
    @(44890, 2022, -0.7009264669202708, 2023, -0.0234350458557242),
    @(45254, 2023, 0.3928252664241905, 2024, 0.196134499498668),
    @(45618, 2024, 0.3224026462283813, 2025, -0.6671574593505647)
)

$startRow = 2

# Row 19 is brand new - give column A the same style (date number format,
# bold, border, centered alignment) already used by A2:A18 before filling it in.
$ws.Range("A18").Copy() | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    if ($null -ne $row[4]) {
        $ws.Cells.Item($r, 5).Value = $row[4]
    } else {
        $ws.Cells.Item($r, 5).ClearContents()
    }
}
